$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 11.651608
$ws.Range("N2").Value = 34.954824
$ws.Range("O2").Value = 0.1892813629236475
$ws.Range("P2").Value = 0.1892813629236474
$ws.Range("Q2").Value = 0.4066022805066667
$ws.Range("R2").Value = 3.65942052456
$ws.Range("S2").Value = 0.1892813629236475
$ws.Range("T2").Value = 0.1892813629236474

$ws.Range("O3").Value = 0.4419371310876561
$ws.Range("P3").Value = 0.4419371310876561
$ws.Range("S3").Value = 0.4419371310876561
$ws.Range("T3").Value = 0.4419371310876561

$ws.Range("M4").Value = 8.657178999999999
$ws.Range("N4").Value = 25.971537
$ws.Range("O4").Value = 0.1406366091439035
$ws.Range("P4").Value = 0.1406366091439035
$ws.Range("Q4").Value = 0.3021066898366667
$ws.Range("R4").Value = 2.71896020853
$ws.Range("S4").Value = 0.1406366091439035
$ws.Range("T4").Value = 0.1406366091439035

$ws.Range("M5").Value = 5.488499666666667
$ws.Range("N5").Value = 16.465499
$ws.Range("O5").Value = 0.08916114387925267
$ws.Range("P5").Value = 0.08916114387925267
$ws.Range("Q5").Value = 0.1915303433677778
$ws.Range("R5").Value = 1.72377309031
$ws.Range("S5").Value = 0.08916114387925267
$ws.Range("T5").Value = 0.08916114387925267

$ws.Range("M6").Value = 4.091608333333333
$ws.Range("N6").Value = 12.274825
$ws.Range("O6").Value = 0.06646852536431769
$ws.Range("P6").Value = 0.06646852536431769
$ws.Range("Q6").Value = 0.1427834921388889
$ws.Range("R6").Value = 1.28505142925
$ws.Range("S6").Value = 0.06646852536431769
$ws.Range("T6").Value = 0.06646852536431769

$ws.Range("M7").Value = 4.463825666666667
$ws.Range("N7").Value = 13.391477
$ws.Range("O7").Value = 0.07251522760122259
$ws.Range("P7").Value = 0.07251522760122257
$ws.Range("Q7").Value = 0.1557726363477778
$ws.Range("R7").Value = 1.40195372713
$ws.Range("S7").Value = 0.07251522760122259
$ws.Range("T7").Value = 0.07251522760122257
